$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Importe (H) column keeps storing these values as text,
# not auto-converted numbers, since the source values are scraped text.
$ws.Range("H2:H153").NumberFormat = "@"

$updates = @{
    'E26' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'F26' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'E33' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'F33' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'E46' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'F46' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'E64' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'F64' = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
    'E63' = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
    'E65' = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
    'E92' = 'RICCOTTI. MARIANA EDITH'
    'H2' = '1475.00'
    'H3' = '8407.90'
    'H4' = '1460.00'
    'H5' = '24.50'
    'H6' = '4333.16'
    'H7' = '3605.00'
    'H8' = '12421.00'
    'H9' = '6453.90'
    'H10' = '45.00'
    'H11' = '13804.80'
    'H12' = '29540.00'
    'H13' = '134100.28'
    'H14' = '242.66'
    'H15' = '490.40'
    'H16' = '297.00'
    'H17' = '13961.67'
    'H18' = '655.00'
    'H19' = '5580.00'
    'H20' = '495.00'
    'H21' = '20506.80'
    'H22' = '138.00'
    'H23' = '92.00'
    'H24' = '21512.59'
    'H25' = '291.34'
    'H26' = '1437.93'
    'H27' = '327.51'
    'H28' = '73420.00'
    'H29' = '5004.00'
    'H30' = '5325.00'
    'H31' = '734.00'
    'H32' = '2894.00'
    'H33' = '456.66'
    'H34' = '38616.83'
    'H35' = '194.00'
    'H36' = '1857.60'
    'H37' = '1884.50'
    'H38' = '11.50'
    'H39' = '1250.00'
    'H40' = '7611.80'
    'H41' = '922.00'
    'H42' = '66.54'
    'H43' = '199.00'
    'H44' = '2347.00'
    'H45' = '9.50'
    'H46' = '8.76'
    'H47' = '25.08'
    'H48' = '12606.40'
    'H49' = '200.00'
    'H50' = '5220.00'
    'H51' = '99.39'
    'H52' = '427.00'
    'H53' = '540.00'
    'H54' = '74.40'
    'H55' = '744.85'
    'H56' = '7540.00'
    'H57' = '140.00'
    'H58' = '145.00'
    'H59' = '400.00'
    'H60' = '19320.00'
    'H61' = '14413.00'
    'H62' = '348.00'
    'H63' = '7198.00'
    'H64' = '154.46'
    'H65' = '1765.00'
    'H66' = '250.00'
    'H67' = '960.00'
    'H68' = '7.53'
    'H69' = '0.02'
    'H70' = '107070.00'
    'H71' = '425.00'
    'H72' = '263.55'
    'H73' = '26.57'
    'H74' = '434.50'
    'H75' = '5732.90'
    'H76' = '2901.00'
    'H77' = '9652.00'
    'H78' = '1011.20'
    'H79' = '13314.40'
    'H80' = '1838.20'
    'H81' = '117.79'
    'H82' = '1599.00'
    'H83' = '471.00'
    'H84' = '200.00'
    'H85' = '375.00'
    'H86' = '4500.00'
    'H87' = '1500.00'
    'H88' = '1875.00'
    'H89' = '600.00'
    'H90' = '1400.00'
    'H91' = '14700.00'
    'H92' = '10000.00'
    'H93' = '4300.00'
    'H94' = '2650.00'
    'H95' = '85.69'
    'H96' = '600.60'
    'H97' = '2.90'
    'H98' = '1568.00'
    'H99' = '56658.00'
    'H100' = '1200.00'
    'H101' = '250.00'
    'H102' = '1000.00'
    'H103' = '1500.00'
    'H104' = '1326.00'
    'H105' = '750.00'
    'H106' = '8000.00'
    'H107' = '6152.66'
    'H108' = '1910.00'
    'H109' = '1200.00'
    'H110' = '400.00'
    'H111' = '1870.00'
    'H112' = '250.00'
    'H113' = '120.00'
    'H114' = '5015.00'
    'H115' = '267.00'
    'H116' = '3144.00'
    'H117' = '350.00'
    'H118' = '14768.48'
    'H119' = '27.00'
    'H120' = '871.58'
    'H121' = '1880.00'
    'H122' = '146.00'
    'H123' = '4174.00'
    'H124' = '32.68'
    'H125' = '5.78'
    'H126' = '527.00'
    'H127' = '4283.55'
    'H128' = '11844.00'
    'H129' = '115.20'
    'H130' = '2861.51'
    'H131' = '1818.63'
    'H132' = '2500.00'
    'H133' = '380000.00'
    'H134' = '69000.00'
    'H135' = '197625.00'
    'H136' = '34000.00'
    'H137' = '346864.00'
    'H138' = '69000.00'
    'H139' = '69000.00'
    'H140' = '146683.50'
    'H141' = '302150.00'
    'H142' = '268416.00'
    'H143' = '380000.00'
    'H144' = '360000.00'
    'H145' = '224444.00'
    'H146' = '162000.00'
    'H147' = '878.26'
    'H148' = '6000.00'
    'H149' = '3672.00'
    'H150' = '6000.00'
    'H151' = '50.00'
    'H152' = '617.00'
    'H153' = '1500.00'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
